$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.074518
$ws.Range("H2").Value = 45.22355399999999
$ws.Range("I2").Value = 0.1146175147123412
$ws.Range("J2").Value = 0.1146175147123412
$ws.Range("M2").Value = 0.118976
$ws.Range("N2").Value = 0.356928
$ws.Range("O2").Value = 0.1290366429533589
$ws.Range("P2").Value = 0.1290366429533589
$ws.Range("Q2").Value = 1.793505853568
$ws.Range("R2").Value = 16.141552682112
$ws.Range("S2").Value = 0.01478985932213773
$ws.Range("T2").Value = 0.01478985932213773
$ws.Range("G3").Value = 15.074518
$ws.Range("H3").Value = 45.22355399999999
$ws.Range("I3").Value = 0.1146175147123412
$ws.Range("J3").Value = 0.1146175147123412
$ws.Range("O3").Value = 0.7268198740608612
$ws.Range("P3").Value = 0.7268198740608612
$ws.Range("Q3").Value = 10.10221336189666
$ws.Range("R3").Value = 90.91992025706998
$ws.Range("S3").Value = 0.08330628760839275
$ws.Range("T3").Value = 0.08330628760839272
$ws.Range("G4").Value = 15.074518
$ws.Range("H4").Value = 45.22355399999999
$ws.Range("I4").Value = 0.1146175147123412
$ws.Range("J4").Value = 0.1146175147123412
$ws.Range("M4").Value = 0.13281
$ws.Range("N4").Value = 0.39843
$ws.Range("O4").Value = 0.1440404497599145
$ws.Range("P4").Value = 0.1440404497599145
$ws.Range("Q4").Value = 2.00204673558
$ws.Range("R4").Value = 18.01842062022
$ws.Range("S4").Value = 0.01650955836952924
$ws.Range("T4").Value = 0.01650955836952924
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("G5").Value = 15.074518
$ws.Range("H5").Value = 45.22355399999999
$ws.Range("I5").Value = 0.1146175147123412
$ws.Range("J5").Value = 0.1146175147123412
$ws.Range("M5").Value = [double]"9.499999999999999E-05"
$ws.Range("N5").Value = 0.000285
$ws.Range("O5").Value = 0.000103033225865461
$ws.Range("P5").Value = 0.000103033225865461
$ws.Range("Q5").Value = 0.00143207921
$ws.Range("R5").Value = 0.01288871289
$ws.Range("S5").Value = [double]"1.180941228149445E-05"
$ws.Range("T5").Value = [double]"1.180941228149445E-05"
$ws.Range("I6").Value = 0.2151681023954678
$ws.Range("J6").Value = 0.2151681023954678
$ws.Range("M6").Value = 0.118976
$ws.Range("N6").Value = 0.356928
$ws.Range("O6").Value = 0.1290366429533589
$ws.Range("P6").Value = 0.1290366429533589
$ws.Range("Q6").Value = 3.366895994176
$ws.Range("R6").Value = 30.302063947584
$ws.Range("S6").Value = 0.02776456960375574
$ws.Range("T6").Value = 0.02776456960375574
$ws.Range("I7").Value = 0.2151681023954678
$ws.Range("J7").Value = 0.2151681023954678
$ws.Range("O7").Value = 0.7268198740608612
$ws.Range("P7").Value = 0.7268198740608612
$ws.Range("S7").Value = 0.1563884530849884
$ws.Range("T7").Value = 0.1563884530849884
$ws.Range("I8").Value = 0.2151681023954678
$ws.Range("J8").Value = 0.2151681023954678
$ws.Range("M8").Value = 0.13281
$ws.Range("N8").Value = 0.39843
$ws.Range("O8").Value = 0.1440404497599145
$ws.Range("P8").Value = 0.1440404497599145
$ws.Range("Q8").Value = 3.758383682309999
$ws.Range("R8").Value = 33.82545314079
$ws.Range("S8").Value = 0.03099291024303052
$ws.Range("T8").Value = 0.03099291024303052
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("I9").Value = 0.2151681023954678
$ws.Range("J9").Value = 0.2151681023954678
$ws.Range("M9").Value = [double]"9.499999999999999E-05"
$ws.Range("N9").Value = 0.000285
$ws.Range("O9").Value = 0.000103033225865461
$ws.Range("P9").Value = 0.000103033225865461
$ws.Range("Q9").Value = 0.002688400345
$ws.Range("R9").Value = 0.024195603105
$ws.Range("S9").Value = [double]"2.216946369315488E-05"
$ws.Range("T9").Value = [double]"2.216946369315488E-05"
$ws.Range("G10").Value = 31.40962533333334
$ws.Range("H10").Value = 94.22887600000001
$ws.Range("I10").Value = 0.2388197880524246
$ws.Range("J10").Value = 0.2388197880524245
$ws.Range("M10").Value = 0.118976
$ws.Range("N10").Value = 0.356928
$ws.Range("O10").Value = 0.1290366429533589
$ws.Range("P10").Value = 0.1290366429533589
$ws.Range("Q10").Value = 3.736991583658668
$ws.Range("R10").Value = 33.63292425292801
$ws.Range("S10").Value = 0.03081650372111755
$ws.Range("T10").Value = 0.03081650372111755
$ws.Range("G11").Value = 31.40962533333334
$ws.Range("H11").Value = 94.22887600000001
$ws.Range("I11").Value = 0.2388197880524246
$ws.Range("J11").Value = 0.2388197880524245
$ws.Range("O11").Value = 0.7268198740608612
$ws.Range("P11").Value = 0.7268198740608612
$ws.Range("Q11").Value = 21.04921276650889
$ws.Range("R11").Value = 189.44291489858
$ws.Range("S11").Value = 0.1735789682755048
$ws.Range("T11").Value = 0.1735789682755048
$ws.Range("G12").Value = 31.40962533333334
$ws.Range("H12").Value = 94.22887600000001
$ws.Range("I12").Value = 0.2388197880524246
$ws.Range("J12").Value = 0.2388197880524245
$ws.Range("M12").Value = 0.13281
$ws.Range("N12").Value = 0.39843
$ws.Range("O12").Value = 0.1440404497599145
$ws.Range("P12").Value = 0.1440404497599145
$ws.Range("Q12").Value = 4.171512340520001
$ws.Range("R12").Value = 37.54361106468
$ws.Range("S12").Value = 0.03439970968263869
$ws.Range("T12").Value = 0.03439970968263869
$ws.Range("D13").Value = "Neutrophils"
$ws.Range("G13").Value = 31.40962533333334
$ws.Range("H13").Value = 94.22887600000001
$ws.Range("I13").Value = 0.2388197880524246
$ws.Range("J13").Value = 0.2388197880524245
$ws.Range("M13").Value = [double]"9.499999999999999E-05"
$ws.Range("N13").Value = 0.000285
$ws.Range("O13").Value = 0.000103033225865461
$ws.Range("P13").Value = 0.000103033225865461
$ws.Range("Q13").Value = 0.002983914406666667
$ws.Range("R13").Value = 0.02685522966
$ws.Range("S13").Value = [double]"2.460637316354699E-05"
$ws.Range("T13").Value = [double]"2.460637316354699E-05"
$ws.Range("G14").Value = 8.433252666666666
$ws.Range("H14").Value = 25.299758
$ws.Range("I14").Value = 0.06412135111680237
$ws.Range("J14").Value = 0.06412135111680237
$ws.Range("M14").Value = 0.118976
$ws.Range("N14").Value = 0.356928
$ws.Range("O14").Value = 0.1290366429533589
$ws.Range("P14").Value = 0.1290366429533589
$ws.Range("Q14").Value = 1.003354669269333
$ws.Range("R14").Value = 9.030192023424
$ws.Range("S14").Value = 0.008274003889745786
$ws.Range("T14").Value = 0.008274003889745786
$ws.Range("G15").Value = 8.433252666666666
$ws.Range("H15").Value = 25.299758
$ws.Range("I15").Value = 0.06412135111680237
$ws.Range("J15").Value = 0.06412135111680237
$ws.Range("O15").Value = 0.7268198740608612
$ws.Range("P15").Value = 0.7268198740608612
$ws.Range("Q15").Value = 5.651558329987777
$ws.Range("R15").Value = 50.86402496989
$ws.Range("S15").Value = 0.04660467234332656
$ws.Range("T15").Value = 0.04660467234332656
$ws.Range("G16").Value = 8.433252666666666
$ws.Range("H16").Value = 25.299758
$ws.Range("I16").Value = 0.06412135111680237
$ws.Range("J16").Value = 0.06412135111680237
$ws.Range("M16").Value = 0.13281
$ws.Range("N16").Value = 0.39843
$ws.Range("O16").Value = 0.1440404497599145
$ws.Range("P16").Value = 0.1440404497599145
$ws.Range("Q16").Value = 1.12002028666
$ws.Range("R16").Value = 10.08018257994
$ws.Range("S16").Value = 0.00923606825407761
$ws.Range("T16").Value = 0.00923606825407761
$ws.Range("D17").Value = "Neutrophils"
$ws.Range("G17").Value = 8.433252666666666
$ws.Range("H17").Value = 25.299758
$ws.Range("I17").Value = 0.06412135111680237
$ws.Range("J17").Value = 0.06412135111680237
$ws.Range("M17").Value = [double]"9.499999999999999E-05"
$ws.Range("N17").Value = 0.000285
$ws.Range("O17").Value = 0.000103033225865461
$ws.Range("P17").Value = 0.000103033225865461
$ws.Range("Q17").Value = 0.0008011590033333332
$ws.Range("R17").Value = 0.00721043103
$ws.Range("S17").Value = [double]"6.60662965241603E-06"
$ws.Range("T17").Value = [double]"6.60662965241603E-06"
$ws.Range("G18").Value = 15.61063466666667
$ws.Range("H18").Value = 46.831904
$ws.Range("I18").Value = 0.1186938214923788
$ws.Range("J18").Value = 0.1186938214923787
$ws.Range("M18").Value = 0.118976
$ws.Range("N18").Value = 0.356928
$ws.Range("O18").Value = 0.1290366429533589
$ws.Range("P18").Value = 0.1290366429533589
$ws.Range("Q18").Value = 1.857290870101334
$ws.Range("R18").Value = 16.715617830912
$ws.Range("S18").Value = 0.01531585226468179
$ws.Range("T18").Value = 0.01531585226468179
$ws.Range("G19").Value = 15.61063466666667
$ws.Range("H19").Value = 46.831904
$ws.Range("I19").Value = 0.1186938214923788
$ws.Range("J19").Value = 0.1186938214923787
$ws.Range("O19").Value = 0.7268198740608612
$ws.Range("P19").Value = 0.7268198740608612
$ws.Range("Q19").Value = 10.46149283959111
$ws.Range("R19").Value = 94.15343555631999
$ws.Range("S19").Value = 0.08626902838889308
$ws.Range("T19").Value = 0.08626902838889305
$ws.Range("G20").Value = 15.61063466666667
$ws.Range("H20").Value = 46.831904
$ws.Range("I20").Value = 0.1186938214923788
$ws.Range("J20").Value = 0.1186938214923787
$ws.Range("M20").Value = 0.13281
$ws.Range("N20").Value = 0.39843
$ws.Range("O20").Value = 0.1440404497599145
$ws.Range("P20").Value = 0.1440404497599145
$ws.Range("Q20").Value = 2.07324839008
$ws.Range("R20").Value = 18.65923551072
$ws.Range("S20").Value = 0.01709671143148524
$ws.Range("T20").Value = 0.01709671143148524
$ws.Range("D21").Value = "Neutrophils"
$ws.Range("G21").Value = 15.61063466666667
$ws.Range("H21").Value = 46.831904
$ws.Range("I21").Value = 0.1186938214923788
$ws.Range("J21").Value = 0.1186938214923787
$ws.Range("M21").Value = [double]"9.499999999999999E-05"
$ws.Range("N21").Value = 0.000285
$ws.Range("O21").Value = 0.000103033225865461
$ws.Range("P21").Value = 0.000103033225865461
$ws.Range("Q21").Value = 0.001483010293333333
$ws.Range("R21").Value = 0.01334709264
$ws.Range("S21").Value = [double]"1.222940731865897E-05"
$ws.Range("T21").Value = [double]"1.222940731865897E-05"
$ws.Range("G22").Value = 32.69321433333334
$ws.Range("H22").Value = 98.079643
$ws.Range("I22").Value = 0.2485794222305853
$ws.Range("J22").Value = 0.2485794222305853
$ws.Range("M22").Value = 0.118976
$ws.Range("N22").Value = 0.356928
$ws.Range("O22").Value = 0.1290366429533589
$ws.Range("P22").Value = 0.1290366429533589
$ws.Range("Q22").Value = 3.889707868522668
$ws.Range("R22").Value = 35.007370816704
$ws.Range("S22").Value = 0.03207585415192028
$ws.Range("T22").Value = 0.03207585415192027
$ws.Range("G23").Value = 32.69321433333334
$ws.Range("H23").Value = 98.079643
$ws.Range("I23").Value = 0.2485794222305853
$ws.Range("J23").Value = 0.2485794222305853
$ws.Range("O23").Value = 0.7268198740608612
$ws.Range("P23").Value = 0.7268198740608612
$ws.Range("Q23").Value = 21.90941207417389
$ws.Range("R23").Value = 197.184708667565
$ws.Range("S23").Value = 0.1806724643597557
$ws.Range("T23").Value = 0.1806724643597556
$ws.Range("G24").Value = 32.69321433333334
$ws.Range("H24").Value = 98.079643
$ws.Range("I24").Value = 0.2485794222305853
$ws.Range("J24").Value = 0.2485794222305853
$ws.Range("M24").Value = 0.13281
$ws.Range("N24").Value = 0.39843
$ws.Range("O24").Value = 0.1440404497599145
$ws.Range("P24").Value = 0.1440404497599145
$ws.Range("Q24").Value = 4.34198579561
$ws.Range("R24").Value = 39.07787216049
$ws.Range("S24").Value = 0.0358054917791532
$ws.Range("T24").Value = 0.0358054917791532
$ws.Range("D25").Value = "Neutrophils"
$ws.Range("G25").Value = 32.69321433333334
$ws.Range("H25").Value = 98.079643
$ws.Range("I25").Value = 0.2485794222305853
$ws.Range("J25").Value = 0.2485794222305853
$ws.Range("M25").Value = [double]"9.499999999999999E-05"
$ws.Range("N25").Value = 0.000285
$ws.Range("O25").Value = 0.000103033225865461
$ws.Range("P25").Value = 0.000103033225865461
$ws.Range("Q25").Value = 0.003105855361666667
$ws.Range("R25").Value = 0.027952698255
$ws.Range("S25").Value = [double]"2.56119397561897E-05"
$ws.Range("T25").Value = [double]"2.56119397561897E-05"

Write-Output "Applied 272 cell updates"
